$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.627.23"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.888.45"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.85"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4830"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2862"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06544"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.859.00"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07444"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.63"
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.084"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.63"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6655"
$ws.Range("E15").Value = "  +3.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.582.86"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.20"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007579"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.30"
$ws.Range("E20").Value = "  +2.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.090.97"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.267"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.187"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.400"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.15"
$ws.Range("E26").Value = "  +2.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.61"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.953"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("E29").Value = "  +11.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.397"
$ws.Range("E30").Value = "  -2.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.330"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.013"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05061"
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.200"
$ws.Range("E34").Value = "  +5.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7495"
$ws.Range("E35").Value = "  +3.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9997"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01884"
$ws.Range("E38").Value = "  +2.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.646"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9211"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.055"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "107.08"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4268"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.635"
$ws.Range("E45").Value = "  -4.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.417"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "64.40"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1274"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.482"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.931"
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.04"
$ws.Range("E51").Value = "  +0.72%  "
